$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab
$ws.Name = "13SEP24_SWEEP_ON_FUEL"

# Update the raw data folder path
$ws.Range("A2").Value = "D:\HN\AUG24Onward\Thesis-03Dec24\SMPS-HN-Desktop\Data\Raw"

# Update header labels
$ws.Range("A1").Value = "f_add_raw"
$ws.Range("B1").Value = "f_name_raw"

# Move the active selection to reflect latest user interaction
$ws.Range("H12").Select()
